$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 114, pushing the existing row 114 (Black Amber,
# Fecha 44592) down to row 116.
$ws.Rows.Item(114).Resize(2).Insert()

# Row 114: Ciruela - Angeleno - Primera (new weekly entry)
$ws.Cells.Item(114, 1).Value = 5
$ws.Cells.Item(114, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(114, 3).Value = 'Maule'
$ws.Cells.Item(114, 4).Value = 44656
$ws.Cells.Item(114, 5).Value = 7
$ws.Cells.Item(114, 6).Value = 'Fruta'
$ws.Cells.Item(114, 7).Value = 100103
$ws.Cells.Item(114, 8).Value = 'Frutos de hueso (carozo)'
$ws.Cells.Item(114, 9).Value = 100103002
$ws.Cells.Item(114, 10).Value = 'Ciruela'
$ws.Cells.Item(114, 11).Value = 'Angeleno'
$ws.Cells.Item(114, 12).Value = 'Primera'
$ws.Cells.Item(114, 13).Value = 180
$ws.Cells.Item(114, 14).Value = 8000
$ws.Cells.Item(114, 15).Value = 8000
$ws.Cells.Item(114, 16).Value = 8000
$ws.Cells.Item(114, 17).Value = '$/bandeja 18 kilos granel'
$ws.Cells.Item(114, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(114, 19).Value = 444
$ws.Cells.Item(114, 20).Value = 18

# Row 115: Ciruela - Angeleno - Segunda (new weekly entry)
$ws.Cells.Item(115, 1).Value = 5
$ws.Cells.Item(115, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(115, 3).Value = 'Maule'
$ws.Cells.Item(115, 4).Value = 44656
$ws.Cells.Item(115, 5).Value = 7
$ws.Cells.Item(115, 6).Value = 'Fruta'
$ws.Cells.Item(115, 7).Value = 100103
$ws.Cells.Item(115, 8).Value = 'Frutos de hueso (carozo)'
$ws.Cells.Item(115, 9).Value = 100103002
$ws.Cells.Item(115, 10).Value = 'Ciruela'
$ws.Cells.Item(115, 11).Value = 'Angeleno'
$ws.Cells.Item(115, 12).Value = 'Segunda'
$ws.Cells.Item(115, 13).Value = 100
$ws.Cells.Item(115, 14).Value = 6000
$ws.Cells.Item(115, 15).Value = 6000
$ws.Cells.Item(115, 16).Value = 6000
$ws.Cells.Item(115, 17).Value = '$/bandeja 18 kilos granel'
$ws.Cells.Item(115, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(115, 19).Value = 333
$ws.Cells.Item(115, 20).Value = 18
